$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "geben",
    "liegen",
    "hauen",
    "schlafen",
    "spüren",
    "zielen",
    "mauern",
    "schlucken",
    "nähen",
    "schätzen",
    "treiben",
    "ehren",
    "wandern",
    "treten",
    "tropfen",
    "drücken",
    "spielen",
    "achten",
    "freuen",
    "werfen",
    "scheitern",
    "kürzen",
    "sorgen",
    "formen",
    "stellen",
    "tollen",
    "platzen",
    "brauchen",
    "teilen",
    "regnen",
    "erben",
    "rasen",
    "münzen",
    "fangen",
    "grüßen",
    "weichen",
    "retten",
    "feiern",
    "quälen",
    "arten",
    "äußern",
    "pflegen",
    "stehlen",
    "pissen",
    "danken",
    "garen",
    "graben",
    "heißen",
    "stecken",
    "fällen",
    "dringen",
    "rauchen",
    "bergen",
    "spannen",
    "flehen",
    "hören",
    "bluten",
    "landen",
    "holen",
    "kümmern",
    "saufen",
    "zögern",
    "handeln",
    "hoffen",
    "wundern",
    "fischen",
    "wiegen",
    "wachsen",
    "führen",
    "melden",
    "schulden",
    "hassen",
    "jagen",
    "rufen",
    "heben",
    "stammen",
    "wehtun",
    "lenken",
    "lockern",
    "ändern",
    "liefern",
    "flüchten",
    "helfen",
    "herrschen",
    "stürmen",
    "schenken",
    "scheinen",
    "ächzen",
    "kosten",
    "filmen",
    "heulen",
    "schreiten",
    "füttern",
    "backen",
    "streichen",
    "schwören",
    "schwimmen",
    "pfeifen",
    "bieten",
    "deuten",
    "jubeln",
    "bremsen",
    "bellen",
    "schwächen",
    "ärgern",
    "klettern",
    "eignen",
    "zünden",
    "altern",
    "loben",
    "irren",
    "sprengen",
    "lügen",
    "starren",
    "mühen",
    "klingeln",
    "biegen",
    "kichern",
    "posten",
    "wüten",
    "mögen",
    "schaden",
    "klingen",
    "foltern",
    "sterben",
    "bitten",
    "pflanzen",
    "runden",
    "schleppen",
    "schmecken",
    "tauchen",
    "suchen",
    "steuern",
    "segeln",
    "heilen",
    "spinnen",
    "greifen",
    "töten",
    "öffnen",
    "dienen",
    "planen",
    "machen",
    "wirken",
    "knarren",
    "kämpfen",
    "schwingen",
    "räumen",
    "wenden",
    "fehlen",
    "malen",
    "fahren",
    "sprechen",
    "zeigen",
    "boxen",
    "wetten",
    "werden",
    "stoßen",
    "trennen",
    "gelten",
    "bauen",
    "lieben",
    "morden",
    "orten",
    "lesen",
    "drehen",
    "rauben",
    "trauen",
    "kehren",
    "enden",
    "sitzen",
    "folgen",
    "reizen",
    "seufzen",
    "sinken",
    "gründen",
    "ahnen",
    "fließen",
    "klagen",
    "schrecken",
    "schicken",
    "stillen",
    "sichern",
    "decken",
    "schnellen",
    "sperren",
    "fallen",
    "haben",
    "warnen",
    "betteln",
    "tanzen",
    "siegen",
    "treffen"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}
